# Fix a missing closing parenthesis in the "matchTypes(...)" code line on the
# "Example: Constraint Checking" slide (slide 21), inside the "Rectangle 3"
# body placeholder.
#
# Before: "        if (!matchTypes(variable.type(), expr)"
# After : "        if (!matchTypes(variable.type(), expr))"
#
# The run that used to hold "(), expr)" is split in two, matching how the
# author's edit landed in the OOXML: "(), " stays in the original run, and a
# new run "expr))" is appended right after it.

$p = $ppt.ActivePresentation

# Locate the slide that contains the "matchTypes(" code sample (robust to the
# slide being anything other than #21 in some other copy of the deck).
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $t = $shape.TextFrame.TextRange.Text
            if ($t -ne $null -and $t.IndexOf("(), expr)") -ge 0) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$full = $tr.Text

# Locate the broken snippet "(), expr)" (missing the closing paren for the
# outer "matchTypes(" call).
$idx0 = $full.IndexOf("(), expr)")
$start1 = $idx0 + 1

# Split the run "(), expr)" into "(), " (kept as-is) and "expr)" -> "expr))"
# (gains the missing paren), preserving the Consolas font formatting.
$firstPart = $tr.Characters($start1, 4)
$firstPart.Text = "(), "

$secondPart = $tr.Characters($start1 + 4, 5)
$secondPart.Text = "expr))"
